$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header row (row 1): number | usernames | password ---
$ws.Range("B1").Value = "usernames"
$ws.Range("C1").Value = "password"
$ws.Range("A1").Value = "number"
$ws.Range("A1:C1").Font.Bold = $true

# --- New column A: row numbers 1,2,3 (left aligned), blank styled cell on row 5 ---
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A2:A5").HorizontalAlignment = -4131

# --- Row 3 explicit (custom) height, matching row 2/4 default ---
$ws.Rows(3).RowHeight = 14.4

# --- New column B width ---
$ws.Columns("B").ColumnWidth = 17.75

# --- Selection cosmetics ---
[void]$ws.Range("B9").Select()
